$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5854969999999999
$ws.Range("H2").Value = 1.170994
$ws.Range("M2").Value = 17.4294175
$ws.Range("N2").Value = 34.858835
$ws.Range("O2").Value = 0.4529581854295807
$ws.Range("P2").Value = 0.3776014560521451
$ws.Range("Q2").Value = 10.2048716579975
$ws.Range("R2").Value = 40.81948663199
$ws.Range("S2").Value = 0.4529581854295807
$ws.Range("T2").Value = 0.3776014560521451

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5854969999999999
$ws.Range("H3").Value = 1.170994
$ws.Range("O3").Value = 0.1017429801035258
$ws.Range("P3").Value = 0.127224648983019
$ws.Range("Q3").Value = 2.292207288569
$ws.Range("R3").Value = 13.753243731414
$ws.Range("S3").Value = 0.1017429801035258
$ws.Range("T3").Value = 0.127224648983019

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5854969999999999
$ws.Range("H4").Value = 1.170994
$ws.Range("M4").Value = 5.397313
$ws.Range("N4").Value = 16.191939
$ws.Range("O4").Value = 0.1402661392829386
$ws.Range("P4").Value = 0.1753959862028526
$ws.Range("Q4").Value = 3.160110569561
$ws.Range("R4").Value = 18.960663417366
$ws.Range("S4").Value = 0.1402661392829386
$ws.Range("T4").Value = 0.1753959862028526

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5854969999999999
$ws.Range("H5").Value = 1.170994
$ws.Range("M5").Value = 5.69137
$ws.Range("N5").Value = 11.38274
$ws.Range("O5").Value = 0.1479081344978025
$ws.Range("P5").Value = 0.1233012863987851
$ws.Range("Q5").Value = 3.33228006089
$ws.Range("R5").Value = 13.32912024356
$ws.Range("S5").Value = 0.1479081344978025
$ws.Range("T5").Value = 0.1233012863987851

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5854969999999999
$ws.Range("H6").Value = 1.170994
$ws.Range("M6").Value = 1.988496
$ws.Range("N6").Value = 5.965488000000001
$ws.Range("O6").Value = 0.05167731737988258
$ws.Range("P6").Value = 0.06461997237892773
$ws.Range("Q6").Value = 1.164258442512
$ws.Range("R6").Value = 6.985550655072
$ws.Range("S6").Value = 0.05167731737988258
$ws.Range("T6").Value = 0.06461997237892773

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.5854969999999999
$ws.Range("H7").Value = 1.170994
$ws.Range("M7").Value = 4.057513666666667
$ws.Range("N7").Value = 12.172541
$ws.Range("O7").Value = 0.1054472433062699
$ws.Range("P7").Value = 0.1318566499842704
$ws.Range("Q7").Value = 2.375662079292333
$ws.Range("R7").Value = 14.253972475754
$ws.Range("S7").Value = 0.1054472433062699
$ws.Range("T7").Value = 0.1318566499842704
